$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.954.14"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.212.54"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "289.20"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.18"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.510"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.33"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0775"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "2.554.97"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "2.211.83"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "39.877.95"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.52"
$ws.Range("E19").Value = "  +10.67%  "
$ws.Range("D20").Value = "0.0₃0880"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.52"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.26"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.42"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "155.39"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.52"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.91"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0713"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.84"
$ws.Range("E36").Value = "  +7.10%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.74"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  +2.95%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.101.33"
$ws.Range("E41").Value = "  +7.94%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.82"
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.00"
$ws.Range("E44").Value = "  +8.02%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0266"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.42"
$ws.Range("E46").Value = "  +8.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.64"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "2.428.39"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.31"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.88"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.44"
$ws.Range("E51").Value = "  +1.00%  "
